$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.779.30'
$ws.Range('E2').Value = '  -2.01%  '
$ws.Range('D3').Value = '3.426.62'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('E4').Value = '  +0.71%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '548.03'
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '176.51'
$ws.Range('E6').Value = '  -3.99%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.631'
$ws.Range('E7').Value = '  +3.61%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.627'
$ws.Range('E9').Value = '  -1.13%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.150'
$ws.Range('E10').Value = '  +2.71%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '53.32'
$ws.Range('E11').Value = '  -4.17%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0000268'
$ws.Range('E12').Value = '  -1.46%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '9.14'
$ws.Range('E13').Value = '  -1.97%  '
$ws.Range('D14').Value = '3.996.13'
$ws.Range('E14').Value = '  +2.02%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '3.444.21'
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.121'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '18.38'
$ws.Range('E17').Value = '  +1.30%  '
$ws.Range('D18').Value = '65.083.31'
$ws.Range('E18').Value = '  -1.40%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.81'
$ws.Range('E19').Value = '  +1.51%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.979'
$ws.Range('E20').Value = '  -1.35%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '412.43'
$ws.Range('E21').Value = '  +1.34%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.01'
$ws.Range('E22').Value = '  +4.60%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '86.26'
$ws.Range('E23').Value = '  +2.65%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '4.25'
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.48'
$ws.Range('E25').Value = '  +7.44%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '10.70'
$ws.Range('E26').Value = '  -10.76%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.83'
$ws.Range('E27').Value = '  -1.53%  '
$ws.Range('E28').Value = '  -2.33%  '
$ws.Range('E29').Value = '  +5.79%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '29.87'
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.49'
$ws.Range('E31').Value = '  -4.55%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '609.07'
$ws.Range('E32').Value = '  -8.18%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '11.62'
$ws.Range('E33').Value = '  +0.82%  '
$ws.Range('E34').Value = '  -0.34%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '58.83'
$ws.Range('E35').Value = '  +0.89%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '37.09'
$ws.Range('E37').Value = '  -3.17%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.143'
$ws.Range('E38').Value = '  +9.34%  '
$ws.Range('D39').Value = '0.0₃0779'
$ws.Range('E39').Value = '  -4.36%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.376'
$ws.Range('E40').Value = '  -5.52%  '
$ws.Range('D41').Value = '3.241.69'
$ws.Range('E41').Value = '  +7.25%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.01'
$ws.Range('E42').Value = '  +0.91%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.28'
$ws.Range('E43').Value = '  -0.92%  '
$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.79'
$ws.Range('E44').Value = '  -1.81%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.52'
$ws.Range('E45').Value = '  -10.30%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.24'
$ws.Range('E46').Value = '  +0.80%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0410'
$ws.Range('E47').Value = '  -1.35%  '
$ws.Range('E48').Value = '  -1.43%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.132'
$ws.Range('E49').Value = '  +1.75%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '136.88'
$ws.Range('E50').Value = '  -2.36%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '8.33'
$ws.Range('E51').Value = '  -3.43%  '
